$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - bold/bordered style matching existing header cells (copy style from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-24
$data = @{
    2  = @(10, 10)
    3  = @(7, 7)
    4  = @(5, 6)
    5  = @(5, 6)
    6  = @(6, 6)
    7  = @(8, 9)
    8  = @(5, 6)
    9  = @(7, 7)
    10 = @(8, 8)
    11 = @(7, 7)
    12 = @(7, 7)
    13 = @(8, 8)
    14 = @(10, 10)
    15 = @(7, 7)
    16 = @(8, 8)
    17 = @(8, 9)
    18 = @(7, 7)
    19 = @(6, 6)
    20 = @(4, 5)
    21 = @(8, 8)
    22 = @(5, 5)
    23 = @(1, 2)
    24 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
